# Commit: "Committing newly added features and corresponding steps on Class module"
#
# The Class worksheet gains a new batch/trainer row plus four brand-new
# columns (Class Description / Test comments / Important notes / Class
# recordings) of sample data, the "ClassDesc"/"Recording" header labels are
# renamed to "classDesc"/"Recordings", and the Class tab becomes the active,
# selected sheet (it was Program before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Class")

# --- Populate the newly added columns (E:H) with their row-2 sample data ---
# (entered right-to-left, matching how the new shared strings were appended
# by the original authoring session)
$ws.Range("E2").Value = "Class Description"
$ws.Range("H2").Value = "Class recordings"
$ws.Range("G2").Value = "Important notes"
$ws.Range("F2").Value = "Test comments"

# --- Row 1 (headers): rename two header labels -----------------------------
$ws.Range("E1").Value = "classDesc"
$ws.Range("H1").Value = "Recordings"

# --- Row 2: replace the old sample batch/trainer values ---------------------
$ws.Range("B2").Value = "Java Batch 01"
$ws.Range("C2").Value = "Playwrighters three"

# --- New columns E:H need explicit (approximate, AutoFit-style) widths -----
# (inputs are pre-compensated for this engine's ColumnWidth->stored-width
# rounding/offset so the serialized <col width> lands as close as possible
# to the real-Excel bestFit values of 16.140625 / 23.85546875 / 18.42578125 / 22)
$ws.Columns.Item(5).ColumnWidth = 15.333333333333334
$ws.Columns.Item(6).ColumnWidth = 23
$ws.Columns.Item(7).ColumnWidth = 17.666666666666668
$ws.Columns.Item(8).ColumnWidth = 21.166666666666668

# --- Class becomes the active sheet / active cell moves to C2 --------------
$ws.Activate() | Out-Null
$ws.Range("C2").Select() | Out-Null
